$d = $word.ActiveDocument
$d.Content.Find.Execute("29-11=18", $true, $false, $false, $false, $false, $true, 1, $false, "99-44=55", 2) | Out-Null
$d.Content.Find.Execute("96-24=72", $true, $false, $false, $false, $false, $true, 1, $false, "15+8=23", 2) | Out-Null
$d.Content.Find.Execute("46-24=22", $true, $false, $false, $false, $false, $true, 1, $false, "56-19=37", 2) | Out-Null
$d.Content.Find.Execute("15+33=48", $true, $false, $false, $false, $false, $true, 1, $false, "42+7=49", 2) | Out-Null
$d.Content.Find.Execute("38-29=9", $true, $false, $false, $false, $false, $true, 1, $false, "39-11=28", 2) | Out-Null
$d.Content.Find.Execute("17+19=36", $true, $false, $false, $false, $false, $true, 1, $false, "15+30=45", 2) | Out-Null
$d.Content.Find.Execute("35+46=81", $true, $false, $false, $false, $false, $true, 1, $false, "70+8=78", 2) | Out-Null
$d.Content.Find.Execute("31-6=25", $true, $false, $false, $false, $false, $true, 1, $false, "63-18=45", 2) | Out-Null
$d.Content.Find.Execute("34-16=18", $true, $false, $false, $false, $false, $true, 1, $false, "63-25=38", 2) | Out-Null
$d.Content.Find.Execute("8+2=10", $true, $false, $false, $false, $false, $true, 1, $false, "2+89=91", 2) | Out-Null
$d.Content.Find.Execute("5+76=81", $true, $false, $false, $false, $false, $true, 1, $false, "61-22=39", 2) | Out-Null
$d.Content.Find.Execute("26+6=32", $true, $false, $false, $false, $false, $true, 1, $false, "23+14=37", 2) | Out-Null
$d.Content.Find.Execute("50+22=72", $true, $false, $false, $false, $false, $true, 1, $false, "88-25=63", 2) | Out-Null
$d.Content.Find.Execute("96-78=18", $true, $false, $false, $false, $false, $true, 1, $false, "39+1=40", 2) | Out-Null
$d.Content.Find.Execute("8+28=36", $true, $false, $false, $false, $false, $true, 1, $false, "43+25=68", 2) | Out-Null
$d.Content.Find.Execute("36+62=98", $true, $false, $false, $false, $false, $true, 1, $false, "46-11=35", 2) | Out-Null
$d.Content.Find.Execute("7+89=96", $true, $false, $false, $false, $false, $true, 1, $false, "20+0=20", 2) | Out-Null
$d.Content.Find.Execute("25+37=62", $true, $false, $false, $false, $false, $true, 1, $false, "22+34=56", 2) | Out-Null
$d.Content.Find.Execute("41-33=8", $true, $false, $false, $false, $false, $true, 1, $false, "43-4=39", 2) | Out-Null
$d.Content.Find.Execute("65-9=56", $true, $false, $false, $false, $false, $true, 1, $false, "70-50=20", 2) | Out-Null
$d.Content.Find.Execute("18+36=54", $true, $false, $false, $false, $false, $true, 1, $false, "69-24=45", 2) | Out-Null
$d.Content.Find.Execute("67+12=79", $true, $false, $false, $false, $false, $true, 1, $false, "82-50=32", 2) | Out-Null
$d.Content.Find.Execute("4+9=13", $true, $false, $false, $false, $false, $true, 1, $false, "37+55=92", 2) | Out-Null
$d.Content.Find.Execute("1+72=73", $true, $false, $false, $false, $false, $true, 1, $false, "34+3=37", 2) | Out-Null
$d.Content.Find.Execute("73-51=22", $true, $false, $false, $false, $false, $true, 1, $false, "24-10=14", 2) | Out-Null
$d.Content.Find.Execute("20+23=43", $true, $false, $false, $false, $false, $true, 1, $false, "20+4=24", 2) | Out-Null
$d.Content.Find.Execute("59+14=73", $true, $false, $false, $false, $false, $true, 1, $false, "63-35=28", 2) | Out-Null
$d.Content.Find.Execute("92-77=15", $true, $false, $false, $false, $false, $true, 1, $false, "87-70=17", 2) | Out-Null
$d.Content.Find.Execute("76+1=77", $true, $false, $false, $false, $false, $true, 1, $false, "30-5=25", 2) | Out-Null
$d.Content.Find.Execute("48-14=34", $true, $false, $false, $false, $false, $true, 1, $false, "67-13=54", 2) | Out-Null
$d.Content.Find.Execute("17+54=71", $true, $false, $false, $false, $false, $true, 1, $false, "83+0=83", 2) | Out-Null
$d.Content.Find.Execute("52+47=99", $true, $false, $false, $false, $false, $true, 1, $false, "70-48=22", 2) | Out-Null
$d.Content.Find.Execute("48-38=10", $true, $false, $false, $false, $false, $true, 1, $false, "62+31=93", 2) | Out-Null
$d.Content.Find.Execute("19+50=69", $true, $false, $false, $false, $false, $true, 1, $false, "10+52=62", 2) | Out-Null
$d.Content.Find.Execute("81-70=11", $true, $false, $false, $false, $false, $true, 1, $false, "48+12=60", 2) | Out-Null
$d.Content.Find.Execute("75-61=14", $true, $false, $false, $false, $false, $true, 1, $false, "53-2=51", 2) | Out-Null
$d.Content.Find.Execute("35+49=84", $true, $false, $false, $false, $false, $true, 1, $false, "14-9=5", 2) | Out-Null
$d.Content.Find.Execute("81-59=22", $true, $false, $false, $false, $false, $true, 1, $false, "37+57=94", 2) | Out-Null
$d.Content.Find.Execute("4+27=31", $true, $false, $false, $false, $false, $true, 1, $false, "96-50=46", 2) | Out-Null
$d.Content.Find.Execute("60-55=5", $true, $false, $false, $false, $false, $true, 1, $false, "32+26=58", 2) | Out-Null
$d.Content.Find.Execute("55+28=83", $true, $false, $false, $false, $false, $true, 1, $false, "34-30=4", 2) | Out-Null
$d.Content.Find.Execute("14+60=74", $true, $false, $false, $false, $false, $true, 1, $false, "98-3=95", 2) | Out-Null
$d.Content.Find.Execute("50-5=45", $true, $false, $false, $false, $false, $true, 1, $false, "94-89=5", 2) | Out-Null
$d.Content.Find.Execute("96-58=38", $true, $false, $false, $false, $false, $true, 1, $false, "86+6=92", 2) | Out-Null
$d.Content.Find.Execute("69-40=29", $true, $false, $false, $false, $false, $true, 1, $false, "55+43=98", 2) | Out-Null
$d.Content.Find.Execute("98-43=55", $true, $false, $false, $false, $false, $true, 1, $false, "39+43=82", 2) | Out-Null
$d.Content.Find.Execute("68-64=4", $true, $false, $false, $false, $false, $true, 1, $false, "78-17=61", 2) | Out-Null
$d.Content.Find.Execute("61+27=88", $true, $false, $false, $false, $false, $true, 1, $false, "85-51=34", 2) | Out-Null
$d.Content.Find.Execute("64-18=46", $true, $false, $false, $false, $false, $true, 1, $false, "75+12=87", 2) | Out-Null
$d.Content.Find.Execute("1+5=6", $true, $false, $false, $false, $false, $true, 1, $false, "76-68=8", 2) | Out-Null
$d.Content.Find.Execute("63-27=36", $true, $false, $false, $false, $false, $true, 1, $false, "14+6=20", 2) | Out-Null
$d.Content.Find.Execute("63+5=68", $true, $false, $false, $false, $false, $true, 1, $false, "55+0=55", 2) | Out-Null
$d.Content.Find.Execute("49-46=3", $true, $false, $false, $false, $false, $true, 1, $false, "7+75=82", 2) | Out-Null
$d.Content.Find.Execute("95-34=61", $true, $false, $false, $false, $false, $true, 1, $false, "94-11=83", 2) | Out-Null
$d.Content.Find.Execute("17+3=20", $true, $false, $false, $false, $false, $true, 1, $false, "0+15=15", 2) | Out-Null
$d.Content.Find.Execute("34+62=96", $true, $false, $false, $false, $false, $true, 1, $false, "34-26=8", 2) | Out-Null
$d.Content.Find.Execute("19+24=43", $true, $false, $false, $false, $false, $true, 1, $false, "4+64=68", 2) | Out-Null
$d.Content.Find.Execute("49-18=31", $true, $false, $false, $false, $false, $true, 1, $false, "8+30=38", 2) | Out-Null
$d.Content.Find.Execute("20+65=85", $true, $false, $false, $false, $false, $true, 1, $false, "94-86=8", 2) | Out-Null
$d.Content.Find.Execute("69-59=10", $true, $false, $false, $false, $false, $true, 1, $false, "0+53=53", 2) | Out-Null
$d.Content.Find.Execute("44-15=29", $true, $false, $false, $false, $false, $true, 1, $false, "79+3=82", 2) | Out-Null
$d.Content.Find.Execute("16+35=51", $true, $false, $false, $false, $false, $true, 1, $false, "30+27=57", 2) | Out-Null
$d.Content.Find.Execute("74-11=63", $true, $false, $false, $false, $false, $true, 1, $false, "63-42=21", 2) | Out-Null
$d.Content.Find.Execute("5+4=9", $true, $false, $false, $false, $false, $true, 1, $false, "96-19=77", 2) | Out-Null
$d.Content.Find.Execute("9-8=1", $true, $false, $false, $false, $false, $true, 1, $false, "70-24=46", 2) | Out-Null
$d.Content.Find.Execute("88-78=10", $true, $false, $false, $false, $false, $true, 1, $false, "78-34=44", 2) | Out-Null
$d.Content.Find.Execute("7+36=43", $true, $false, $false, $false, $false, $true, 1, $false, "34+50=84", 2) | Out-Null
$d.Content.Find.Execute("29+24=53", $true, $false, $false, $false, $false, $true, 1, $false, "65-22=43", 2) | Out-Null
$d.Content.Find.Execute("20+27=47", $true, $false, $false, $false, $false, $true, 1, $false, "44+46=90", 2) | Out-Null
$d.Content.Find.Execute("11+10=21", $true, $false, $false, $false, $false, $true, 1, $false, "94-63=31", 2) | Out-Null
$d.Content.Find.Execute("43-20=23", $true, $false, $false, $false, $false, $true, 1, $false, "78-52=26", 2) | Out-Null
$d.Content.Find.Execute("58-38=20", $true, $false, $false, $false, $false, $true, 1, $false, "72-32=40", 2) | Out-Null
$d.Content.Find.Execute("63+13=76", $true, $false, $false, $false, $false, $true, 1, $false, "3-1=2", 2) | Out-Null
$d.Content.Find.Execute("27-20=7", $true, $false, $false, $false, $false, $true, 1, $false, "17+71=88", 2) | Out-Null
$d.Content.Find.Execute("75-25=50", $true, $false, $false, $false, $false, $true, 1, $false, "7+41=48", 2) | Out-Null
$d.Content.Find.Execute("14+4=18", $true, $false, $false, $false, $false, $true, 1, $false, "63-54=9", 2) | Out-Null
$d.Content.Find.Execute("24+17=41", $true, $false, $false, $false, $false, $true, 1, $false, "52+36=88", 2) | Out-Null
$d.Content.Find.Execute("30+7=37", $true, $false, $false, $false, $false, $true, 1, $false, "57-51=6", 2) | Out-Null
$d.Content.Find.Execute("20+35=55", $true, $false, $false, $false, $false, $true, 1, $false, "81+12=93", 2) | Out-Null
$d.Content.Find.Execute("94-83=11", $true, $false, $false, $false, $false, $true, 1, $false, "6-1=5", 2) | Out-Null
$d.Content.Find.Execute("4+40=44", $true, $false, $false, $false, $false, $true, 1, $false, "98-82=16", 2) | Out-Null
$d.Content.Find.Execute("91-3=88", $true, $false, $false, $false, $false, $true, 1, $false, "90-50=40", 2) | Out-Null
$d.Content.Find.Execute("16+75=91", $true, $false, $false, $false, $false, $true, 1, $false, "29-10=19", 2) | Out-Null
$d.Content.Find.Execute("51-19=32", $true, $false, $false, $false, $false, $true, 1, $false, "96-51=45", 2) | Out-Null
$d.Content.Find.Execute("0+13=13", $true, $false, $false, $false, $false, $true, 1, $false, "91-6=85", 2) | Out-Null
$d.Content.Find.Execute("30+56=86", $true, $false, $false, $false, $false, $true, 1, $false, "69+0=69", 2) | Out-Null
$d.Content.Find.Execute("20+59=79", $true, $false, $false, $false, $false, $true, 1, $false, "54-26=28", 2) | Out-Null
$d.Content.Find.Execute("63+36=99", $true, $false, $false, $false, $false, $true, 1, $false, "37-4=33", 2) | Out-Null
$d.Content.Find.Execute("35+32=67", $true, $false, $false, $false, $false, $true, 1, $false, "80-74=6", 2) | Out-Null
$d.Content.Find.Execute("53-19=34", $true, $false, $false, $false, $false, $true, 1, $false, "60-21=39", 2) | Out-Null
$d.Content.Find.Execute("37-29=8", $true, $false, $false, $false, $false, $true, 1, $false, "31+51=82", 2) | Out-Null
$d.Content.Find.Execute("1+49=50", $true, $false, $false, $false, $false, $true, 1, $false, "27-11=16", 2) | Out-Null
$d.Content.Find.Execute("63+1=64", $true, $false, $false, $false, $false, $true, 1, $false, "59-49=10", 2) | Out-Null
$d.Content.Find.Execute("11+42=53", $true, $false, $false, $false, $false, $true, 1, $false, "33-31=2", 2) | Out-Null
$d.Content.Find.Execute("95-85=10", $true, $false, $false, $false, $false, $true, 1, $false, "60-1=59", 2) | Out-Null
$d.Content.Find.Execute("21-17=4", $true, $false, $false, $false, $false, $true, 1, $false, "81-79=2", 2) | Out-Null
$d.Content.Find.Execute("26+1=27", $true, $false, $false, $false, $false, $true, 1, $false, "93-3=90", 2) | Out-Null
$d.Content.Find.Execute("34+59=93", $true, $false, $false, $false, $false, $true, 1, $false, "81-14=67", 2) | Out-Null
$d.Content.Find.Execute("69-27=42", $true, $false, $false, $false, $false, $true, 1, $false, "18+57=75", 2) | Out-Null
$d.Content.Find.Execute("20+22=42", $true, $false, $false, $false, $false, $true, 1, $false, "1+68=69", 2) | Out-Null
